$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2, A3, A4 with combined tuple-style strings
$ws.Range("A2").Value = "(""Gaea's Cradle"", ['Legendary Land', '{T}: Add {G} for each creature you control.'])"
$ws.Range("A3").Value = "('Lightning Bolt', ['{R}', 'Instant', 'Lightning Bolt deals 3 damage to any target.'])"
$ws.Range("A4").Value = "('Stroke of Genius', ['{X}{2}{U}', 'Instant', 'Target player draws X cards.'])"

# Clear rows 5 through 12, which are no longer needed
$ws.Range("A5:A12").ClearContents()
